$d = $word.ActiveDocument

# --- Create the new character styles ---

$charStyle = [Microsoft.Office.Interop.Word.WdStyleType]::wdStyleTypeCharacter

$gaNStyle = $d.Styles.Add("GaNStyle", $charStyle)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", $charStyle)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", $charStyle)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = [Microsoft.Office.Interop.Word.WdUnderline]::wdUnderlineSingle

# --- Apply GaNStyle to every "2022 Fechas de la campaña..." run ---

$rng = $d.Range(0, 0)
$iterations = 0
while ($rng.Find.Execute("2022 Fechas de la campaña para constelación de pegaso: 8-17 de octubre, 7-16 de noviembre,", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $iterations = $iterations + 1
    if ($iterations -gt 20) { break }
}

# --- Apply GaNParagraph to the "Usted está participando..." paragraph run ---

$rng = $d.Range(0, 0)
if ($rng.Find.Execute("Usted está participando en una campaña mundial", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $para = $rng.Paragraphs(1)
    $paraRange = $para.Range
    $target = $d.Range($paraRange.Start, $paraRange.End - 1)
    $target.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the GaNight link run ---

$rng = $d.Range(0, 0)
if ($rng.Find.Execute("(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNLinks"
}
